# Update weekly fruit/vegetable prices (Membrillo - Vega Monumental Concepción)
# The row data (Fecha, Calidad, Volumen, Precio min/max/prom, Unidad, Origen, Precio $/Kg, Kg/unidad)
# has been rotated across rows 2 and 4-12 (row 3 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, columns: D, L, M, N, O, P, Q, R, S, T
$rows = @{
    2  = @(44698, "Primera", 50,  10000, 10000, 10000, "`$/caja 18 kilos granel",      "Región de O'Higgins", 556, 18)
    4  = @(44299, "Primera", 100, 10000, 11000, 10500, "`$/caja 18 kilos granel",      "Región del Maule",    583, 18)
    5  = @(44299, "Segunda", 50,  9000,  9000,  9000,  "`$/caja 18 kilos granel",      "Región del Maule",    500, 18)
    6  = @(44307, "Primera", 50,  10000, 10000, 10000, "`$/bandeja 18 kilos granel",   "Región de O'Higgins", 556, 18)
    7  = @(44307, "Segunda", 50,  8000,  8000,  8000,  "`$/bandeja 18 kilos granel",   "Región de O'Higgins", 444, 18)
    8  = @(44272, "Primera", 100, 9000,  10000, 9500,  "`$/caja 15 kilos granel",      "Región de O'Higgins", 633, 15)
    9  = @(44272, "Segunda", 50,  8000,  8000,  8000,  "`$/caja 15 kilos granel",      "Región de O'Higgins", 533, 15)
    10 = @(44363, "Primera", 100, 9000,  10000, 9500,  "`$/caja 15 kilos empedrada",   "Región de O'Higgins", 633, 15)
    11 = @(44425, "Primera", 100, 12000, 13000, 12500, "`$/bandeja 18 kilos granel",   "Región de O'Higgins", 694, 18)
    12 = @(44358, "Primera", 100, 11000, 12000, 11500, "`$/caja 18 kilos granel",      "Región de O'Higgins", 639, 18)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Range("D$r").Value  = $vals[0]
    $ws.Range("L$r").Value  = $vals[1]
    $ws.Range("M$r").Value  = $vals[2]
    $ws.Range("N$r").Value  = $vals[3]
    $ws.Range("O$r").Value  = $vals[4]
    $ws.Range("P$r").Value  = $vals[5]
    $ws.Range("Q$r").Value  = $vals[6]
    $ws.Range("R$r").Value  = $vals[7]
    $ws.Range("S$r").Value  = $vals[8]
    $ws.Range("T$r").Value  = $vals[9]
}
